$d = $word.ActiveDocument
$d.Content.Find.Execute("60+17=", $true, $false, $false, $false, $false, $true, 1, $false, "89-34=", 2)
$d.Content.Find.Execute("59+20=", $true, $false, $false, $false, $false, $true, 1, $false, "48+28=", 2)
$d.Content.Find.Execute("26+70=", $true, $false, $false, $false, $false, $true, 1, $false, "85-34=", 2)
$d.Content.Find.Execute("64-44=", $true, $false, $false, $false, $false, $true, 1, $false, "74-29=", 2)
$d.Content.Find.Execute("3+54=", $true, $false, $false, $false, $false, $true, 1, $false, "41+44=", 2)
$d.Content.Find.Execute("74-27=", $true, $false, $false, $false, $false, $true, 1, $false, "70-48=", 2)
$d.Content.Find.Execute("69-33=", $true, $false, $false, $false, $false, $true, 1, $false, "40-39=", 2)
$d.Content.Find.Execute("62-32=", $true, $false, $false, $false, $false, $true, 1, $false, "3+48=", 2)
$d.Content.Find.Execute("45-3=", $true, $false, $false, $false, $false, $true, 1, $false, "1+47=", 2)
$d.Content.Find.Execute("14+81=", $true, $false, $false, $false, $false, $true, 1, $false, "0+7=", 2)
$d.Content.Find.Execute("77-38=", $true, $false, $false, $false, $false, $true, 1, $false, "60+7=", 2)
$d.Content.Find.Execute("33+21=", $true, $false, $false, $false, $false, $true, 1, $false, "53+43=", 2)
$d.Content.Find.Execute("35+28=", $true, $false, $false, $false, $false, $true, 1, $false, "92-91=", 2)
$d.Content.Find.Execute("25+11=", $true, $false, $false, $false, $false, $true, 1, $false, "43+15=", 2)
$d.Content.Find.Execute("25+39=", $true, $false, $false, $false, $false, $true, 1, $false, "16-7=", 2)
$d.Content.Find.Execute("94-22=", $true, $false, $false, $false, $false, $true, 1, $false, "81-67=", 2)
$d.Content.Find.Execute("21+21=", $true, $false, $false, $false, $false, $true, 1, $false, "61+31=", 2)
$d.Content.Find.Execute("26+42=", $true, $false, $false, $false, $false, $true, 1, $false, "97-34=", 2)
$d.Content.Find.Execute("94-17=", $true, $false, $false, $false, $false, $true, 1, $false, "41+18=", 2)
$d.Content.Find.Execute("25+45=", $true, $false, $false, $false, $false, $true, 1, $false, "5+61=", 2)
$d.Content.Find.Execute("28+66=", $true, $false, $false, $false, $false, $true, 1, $false, "12+86=", 2)
$d.Content.Find.Execute("87-84=", $true, $false, $false, $false, $false, $true, 1, $false, "61-2=", 2)
$d.Content.Find.Execute("64-54=", $true, $false, $false, $false, $false, $true, 1, $false, "54-4=", 2)
$d.Content.Find.Execute("80-1=", $true, $false, $false, $false, $false, $true, 1, $false, "83-78=", 2)
$d.Content.Find.Execute("52-50=", $true, $false, $false, $false, $false, $true, 1, $false, "10+38=", 2)
$d.Content.Find.Execute("35-2=", $true, $false, $false, $false, $false, $true, 1, $false, "70-2=", 2)
$d.Content.Find.Execute("60-12=", $true, $false, $false, $false, $false, $true, 1, $false, "59+32=", 2)
$d.Content.Find.Execute("68-57=", $true, $false, $false, $false, $false, $true, 1, $false, "30+69=", 2)
$d.Content.Find.Execute("61+5=", $true, $false, $false, $false, $false, $true, 1, $false, "55+27=", 2)
$d.Content.Find.Execute("41-17=", $true, $false, $false, $false, $false, $true, 1, $false, "7+23=", 2)
$d.Content.Find.Execute("36+27=", $true, $false, $false, $false, $false, $true, 1, $false, "82-35=", 2)
$d.Content.Find.Execute("52-27=", $true, $false, $false, $false, $false, $true, 1, $false, "10+30=", 2)
$d.Content.Find.Execute("43-9=", $true, $false, $false, $false, $false, $true, 1, $false, "91-36=", 2)
$d.Content.Find.Execute("44-20=", $true, $false, $false, $false, $false, $true, 1, $false, "50-14=", 2)
$d.Content.Find.Execute("53+42=", $true, $false, $false, $false, $false, $true, 1, $false, "46+31=", 2)
$d.Content.Find.Execute("76-8=", $true, $false, $false, $false, $false, $true, 1, $false, "91-30=", 2)
$d.Content.Find.Execute("56+8=", $true, $false, $false, $false, $false, $true, 1, $false, "79-20=", 2)
$d.Content.Find.Execute("82-73=", $true, $false, $false, $false, $false, $true, 1, $false, "23+22=", 2)
$d.Content.Find.Execute("27+28=", $true, $false, $false, $false, $false, $true, 1, $false, "82-33=", 2)
$d.Content.Find.Execute("75-38=", $true, $false, $false, $false, $false, $true, 1, $false, "93-31=", 2)
$d.Content.Find.Execute("77-10=", $true, $false, $false, $false, $false, $true, 1, $false, "83-49=", 2)
$d.Content.Find.Execute("22+51=", $true, $false, $false, $false, $false, $true, 1, $false, "99-65=", 2)
$d.Content.Find.Execute("30+20=", $true, $false, $false, $false, $false, $true, 1, $false, "25+73=", 2)
$d.Content.Find.Execute("30-3=", $true, $false, $false, $false, $false, $true, 1, $false, "34+64=", 2)
$d.Content.Find.Execute("0+77=", $true, $false, $false, $false, $false, $true, 1, $false, "82-10=", 2)
$d.Content.Find.Execute("42-32=", $true, $false, $false, $false, $false, $true, 1, $false, "28+2=", 2)
$d.Content.Find.Execute("27+39=", $true, $false, $false, $false, $false, $true, 1, $false, "90-20=", 2)
$d.Content.Find.Execute("23+62=", $true, $false, $false, $false, $false, $true, 1, $false, "27+31=", 2)
$d.Content.Find.Execute("69-46=", $true, $false, $false, $false, $false, $true, 1, $false, "46-11=", 2)
$d.Content.Find.Execute("16+14=", $true, $false, $false, $false, $false, $true, 1, $false, "14+62=", 2)
$d.Content.Find.Execute("38+3=", $true, $false, $false, $false, $false, $true, 1, $false, "49+49=", 2)
$d.Content.Find.Execute("15+53=", $true, $false, $false, $false, $false, $true, 1, $false, "63+30=", 2)
$d.Content.Find.Execute("14+46=", $true, $false, $false, $false, $false, $true, 1, $false, "39-8=", 2)
$d.Content.Find.Execute("74+19=", $true, $false, $false, $false, $false, $true, 1, $false, "75-46=", 2)
$d.Content.Find.Execute("33+13=", $true, $false, $false, $false, $false, $true, 1, $false, "80-6=", 2)
$d.Content.Find.Execute("89-80=", $true, $false, $false, $false, $false, $true, 1, $false, "68-14=", 2)
$d.Content.Find.Execute("95-60=", $true, $false, $false, $false, $false, $true, 1, $false, "20+6=", 2)
$d.Content.Find.Execute("1+84=", $true, $false, $false, $false, $false, $true, 1, $false, "97+1=", 2)
$d.Content.Find.Execute("51-23=", $true, $false, $false, $false, $false, $true, 1, $false, "15+56=", 2)
$d.Content.Find.Execute("84-18=", $true, $false, $false, $false, $false, $true, 1, $false, "75-41=", 2)
$d.Content.Find.Execute("92-16=", $true, $false, $false, $false, $false, $true, 1, $false, "80-44=", 2)
$d.Content.Find.Execute("36+56=", $true, $false, $false, $false, $false, $true, 1, $false, "49-16=", 2)
$d.Content.Find.Execute("83-69=", $true, $false, $false, $false, $false, $true, 1, $false, "92-54=", 2)
$d.Content.Find.Execute("43+43=", $true, $false, $false, $false, $false, $true, 1, $false, "47+31=", 2)
$d.Content.Find.Execute("56+28=", $true, $false, $false, $false, $false, $true, 1, $false, "44+44=", 2)
$d.Content.Find.Execute("57-39=", $true, $false, $false, $false, $false, $true, 1, $false, "60-26=", 2)
$d.Content.Find.Execute("49+23=", $true, $false, $false, $false, $false, $true, 1, $false, "87-25=", 2)
$d.Content.Find.Execute("57-25=", $true, $false, $false, $false, $false, $true, 1, $false, "88-15=", 2)
$d.Content.Find.Execute("64+18=", $true, $false, $false, $false, $false, $true, 1, $false, "44-31=", 2)
$d.Content.Find.Execute("72-48=", $true, $false, $false, $false, $false, $true, 1, $false, "63+29=", 2)
$d.Content.Find.Execute("96+1=", $true, $false, $false, $false, $false, $true, 1, $false, "26+1=", 2)
$d.Content.Find.Execute("1+37=", $true, $false, $false, $false, $false, $true, 1, $false, "77-56=", 2)
$d.Content.Find.Execute("13+51=", $true, $false, $false, $false, $false, $true, 1, $false, "93-32=", 2)
$d.Content.Find.Execute("33-12=", $true, $false, $false, $false, $false, $true, 1, $false, "70-26=", 2)
$d.Content.Find.Execute("5+6=", $true, $false, $false, $false, $false, $true, 1, $false, "26+11=", 2)
$d.Content.Find.Execute("81-27=", $true, $false, $false, $false, $false, $true, 1, $false, "40+38=", 2)
$d.Content.Find.Execute("77+8=", $true, $false, $false, $false, $false, $true, 1, $false, "53-2=", 2)
$d.Content.Find.Execute("32+55=", $true, $false, $false, $false, $false, $true, 1, $false, "36+10=", 2)
$d.Content.Find.Execute("90-62=", $true, $false, $false, $false, $false, $true, 1, $false, "33+45=", 2)
$d.Content.Find.Execute("61-59=", $true, $false, $false, $false, $false, $true, 1, $false, "19-14=", 2)
$d.Content.Find.Execute("28+12=", $true, $false, $false, $false, $false, $true, 1, $false, "91-84=", 2)
$d.Content.Find.Execute("67-14=", $true, $false, $false, $false, $false, $true, 1, $false, "2+31=", 2)
$d.Content.Find.Execute("91-25=", $true, $false, $false, $false, $false, $true, 1, $false, "96-21=", 2)
$d.Content.Find.Execute("14+69=", $true, $false, $false, $false, $false, $true, 1, $false, "50+45=", 2)
$d.Content.Find.Execute("9+31=", $true, $false, $false, $false, $false, $true, 1, $false, "13+40=", 2)
$d.Content.Find.Execute("16+58=", $true, $false, $false, $false, $false, $true, 1, $false, "93-30=", 2)
$d.Content.Find.Execute("81-48=", $true, $false, $false, $false, $false, $true, 1, $false, "26+39=", 2)
$d.Content.Find.Execute("4+58=", $true, $false, $false, $false, $false, $true, 1, $false, "77-5=", 2)
$d.Content.Find.Execute("44+36=", $true, $false, $false, $false, $false, $true, 1, $false, "2+71=", 2)
$d.Content.Find.Execute("65-26=", $true, $false, $false, $false, $false, $true, 1, $false, "58+40=", 2)
$d.Content.Find.Execute("10+41=", $true, $false, $false, $false, $false, $true, 1, $false, "19+34=", 2)
$d.Content.Find.Execute("27+35=", $true, $false, $false, $false, $false, $true, 1, $false, "21+3=", 2)
$d.Content.Find.Execute("50+30=", $true, $false, $false, $false, $false, $true, 1, $false, "44+4=", 2)
$d.Content.Find.Execute("85-23=", $true, $false, $false, $false, $false, $true, 1, $false, "11+84=", 2)
$d.Content.Find.Execute("18+39=", $true, $false, $false, $false, $false, $true, 1, $false, "23-22=", 2)
$d.Content.Find.Execute("95-56=", $true, $false, $false, $false, $false, $true, 1, $false, "98-5=", 2)
$d.Content.Find.Execute("17+2=", $true, $false, $false, $false, $false, $true, 1, $false, "64-55=", 2)
$d.Content.Find.Execute("83+1=", $true, $false, $false, $false, $false, $true, 1, $false, "29+1=", 2)
$d.Content.Find.Execute("67-38=", $true, $false, $false, $false, $false, $true, 1, $false, "58+6=", 2)
$d.Content.Find.Execute("67+6=", $true, $false, $false, $false, $false, $true, 1, $false, "16+0=", 2)
